# Insert a new data row for Frambuesa (raspberry) prices.
# The new record is inserted as row 45 on the active sheet, pushing the
# existing rows 45-103 down to 46-104 (dimension grows from A1:T103 to A1:T104).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(45).Insert()

$ws.Cells.Item(45, 1).Value  = 9
$ws.Cells.Item(45, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(45, 3).Value  = "Metropolitana"
$ws.Cells.Item(45, 4).Value  = 44665
$ws.Cells.Item(45, 5).Value  = 13
$ws.Cells.Item(45, 6).Value  = "Fruta"
$ws.Cells.Item(45, 7).Value  = 100101
$ws.Cells.Item(45, 8).Value  = "Berries"
$ws.Cells.Item(45, 9).Value  = 100101004
$ws.Cells.Item(45, 10).Value = "Frambuesa"
$ws.Cells.Item(45, 11).Value = "Sin especificar"
$ws.Cells.Item(45, 12).Value = "Primera"
$ws.Cells.Item(45, 13).Value = 380
$ws.Cells.Item(45, 14).Value = 8000
$ws.Cells.Item(45, 15).Value = 8000
$ws.Cells.Item(45, 16).Value = 8000
$ws.Cells.Item(45, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(45, 18).Value = "Provincia de Linares"
$ws.Cells.Item(45, 19).Value = 4000
$ws.Cells.Item(45, 20).Value = 2
